$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3217.6428
$ws.Range("I6").Value = 397.1
$ws.Range("J6").Value = 10269
$ws.Range("K6").Value = 1191.3
$ws.Range("L6").Value = 30807
$ws.Range("M6").Value = -1079.3
$ws.Range("N6").Value = -31031
$ws.Range("H8").Value = 5957.75
$ws.Range("I8").Value = 5099.1
$ws.Range("J8").Value = 10251
$ws.Range("K8").Value = 15297.3
$ws.Range("L8").Value = 30753
$ws.Range("M8").Value = -15158.3
$ws.Range("N8").Value = -31031
$ws.Range("H129").Value = 23006.49
$ws.Range("I129").Value = 951.63635
$ws.Range("J129").Value = 29745.473
$ws.Range("K129").Value = 2854.90905
$ws.Range("L129").Value = 89236.41900000001
$ws.Range("M129").Value = 2145.09095
$ws.Range("N129").Value = -99236.41900000001
$ws.Range("H138").Value = 3644.8867
$ws.Range("I138").Value = 2026.0312
$ws.Range("J138").Value = 6111.7144
$ws.Range("K138").Value = 6078.0936
$ws.Range("L138").Value = 18335.1432
$ws.Range("M138").Value = -938.0936000000002
$ws.Range("N138").Value = -28615.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1804.5454
$ws.Range("I110").Value = 1377.7778
$ws.Range("J110").Value = 3725
$ws.Range("K110").Value = 1377.7778
$ws.Range("L110").Value = 3725
$ws.Range("M110").Value = 667.2221999999999
$ws.Range("N110").Value = -7815

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 50006
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 50006
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 50006
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -50354
$ws.Range("H118").Value = 37976.668
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 37976.668
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 37976.668
$ws.Range("N118").Value = -41290.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 125.25
$ws.Range("I40").Value = 125.25
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 501
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -432
$ws.Range("N40").ClearContents()
$ws.Range("H68").Value = 1404.5737
$ws.Range("I68").Value = 1080.2693
$ws.Range("J68").Value = 1645.4857
$ws.Range("K68").Value = 3240.8079
$ws.Range("L68").Value = 4936.4571
$ws.Range("M68").Value = -2429.8079
$ws.Range("N68").Value = -6558.4571
$ws.Range("H71").Value = 1404.5737
$ws.Range("I71").Value = 1080.2693
$ws.Range("J71").Value = 1645.4857
$ws.Range("K71").Value = 9722.423699999999
$ws.Range("L71").Value = 14809.3713
$ws.Range("M71").Value = -5666.423699999999
$ws.Range("N71").Value = -22921.3713
$ws.Range("H86").Value = 75
$ws.Range("I86").Value = 50
$ws.Range("J86").Value = 100
$ws.Range("K86").Value = 150
$ws.Range("L86").Value = 300
$ws.Range("M86").Value = 1036
$ws.Range("N86").Value = -2672
$ws.Range("H87").Value = 22210.066
$ws.Range("I87").Value = 970.2
$ws.Range("J87").Value = 32830
$ws.Range("K87").Value = 2910.6
$ws.Range("L87").Value = 98490
$ws.Range("M87").Value = -1662.6
$ws.Range("N87").Value = -100986
$ws.Range("H89").Value = 75
$ws.Range("I89").Value = 50
$ws.Range("J89").Value = 100
$ws.Range("K89").Value = 450
$ws.Range("L89").Value = 900
$ws.Range("M89").Value = 5478
$ws.Range("N89").Value = -12756
$ws.Range("H90").Value = 22210.066
$ws.Range("I90").Value = 970.2
$ws.Range("J90").Value = 32830
$ws.Range("K90").Value = 8731.800000000001
$ws.Range("L90").Value = 295470
$ws.Range("M90").Value = -2491.800000000001
$ws.Range("N90").Value = -307950
$ws.Range("H104").Value = 1581.6666
$ws.Range("I104").Value = 663.3333
$ws.Range("J104").Value = 2500
$ws.Range("K104").Value = 1989.9999
$ws.Range("L104").Value = 7500
$ws.Range("M104").Value = 631.0001
$ws.Range("N104").Value = -12742
$ws.Range("H105").Value = 8112.0835
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 8112.0835
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 24336.2505
$ws.Range("N105").Value = -29578.2505
$ws.Range("H107").Value = 383.7647
$ws.Range("I107").Value = 206.6
$ws.Range("J107").Value = 1133.3077
$ws.Range("K107").Value = 619.8
$ws.Range("L107").Value = 3399.9231
$ws.Range("M107").Value = 1300.2
$ws.Range("N107").Value = -7239.9231
$ws.Range("H108").Value = 758.7143
$ws.Range("I108").Value = 221.83333
$ws.Range("J108").Value = 3980
$ws.Range("K108").Value = 665.49999
$ws.Range("L108").Value = 11940
$ws.Range("M108").Value = 2214.50001
$ws.Range("N108").Value = -17700
$ws.Range("H109").Value = 1518037.8
$ws.Range("I109").Value = 1319.8
$ws.Range("J109").Value = 2781969.2
$ws.Range("K109").Value = 3959.4
$ws.Range("L109").Value = 8345907.600000001
$ws.Range("M109").Value = -2919.4
$ws.Range("N109").Value = -8347987.600000001
$ws.Range("H110").Value = 6029.5713
$ws.Range("I110").Value = 3289
$ws.Range("J110").Value = 8085
$ws.Range("K110").Value = 9867
$ws.Range("L110").Value = 24255
$ws.Range("M110").Value = -5777
$ws.Range("N110").Value = -32435
$ws.Range("H111").Value = 2817.3333
$ws.Range("I111").Value = 2838.5
$ws.Range("J111").Value = 2775
$ws.Range("K111").Value = 8515.5
$ws.Range("L111").Value = 8325
$ws.Range("M111").Value = -5448.5
$ws.Range("N111").Value = -14459
$ws.Range("H112").Value = 5390
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 5390
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 16170
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -18386
$ws.Range("H114").Value = 1693.3846
$ws.Range("I114").Value = 861.4
$ws.Range("J114").Value = 4466.6665
$ws.Range("K114").Value = 2584.2
$ws.Range("L114").Value = 13399.9995
$ws.Range("M114").Value = 669.8000000000002
$ws.Range("N114").Value = -19907.9995
$ws.Range("H116").Value = 575.5714
$ws.Range("I116").Value = 575.5714
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1726.7142
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1715.2858
$ws.Range("H117").Value = 568
$ws.Range("I117").Value = 568
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 1704
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 1738
$ws.Range("N117").ClearContents()
$ws.Range("H118").Value = 1854717.8
$ws.Range("I118").Value = 1998.2
$ws.Range("J118").Value = 2781077.5
$ws.Range("K118").Value = 5994.6
$ws.Range("L118").Value = 8343232.5
$ws.Range("M118").Value = -4751.6
$ws.Range("N118").Value = -8345718.5
$ws.Range("H119").Value = 2874.875
$ws.Range("I119").Value = 2428.4285
$ws.Range("J119").Value = 6000
$ws.Range("K119").Value = 7285.2855
$ws.Range("L119").Value = 18000
$ws.Range("M119").Value = -2447.2855
$ws.Range("N119").Value = -27676
$ws.Range("H120").Value = 10832.857
$ws.Range("I120").Value = 1276.6666
$ws.Range("J120").Value = 18000
$ws.Range("K120").Value = 3829.9998
$ws.Range("L120").Value = 54000
$ws.Range("M120").Value = 1008.0002
$ws.Range("N120").Value = -63676
$ws.Range("H121").Value = 31250832
$ws.Range("I121").Value = 387.5
$ws.Range("J121").Value = 35715184
$ws.Range("K121").Value = 1162.5
$ws.Range("L121").Value = 107145552
$ws.Range("M121").Value = 147.5
$ws.Range("N121").Value = -107148172
$ws.Range("H124").Value = 1044.2858
$ws.Range("I124").Value = 868.3333
$ws.Range("J124").Value = 2100
$ws.Range("K124").Value = 2604.9999
$ws.Range("L124").Value = 6300
$ws.Range("M124").Value = 2305.0001
$ws.Range("N124").Value = -16120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 46500
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 46500
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 46500
$ws.Range("N133").Value = -56620
